$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.299.07'
$ws.Range('E2').Value = '  -7.63%  '
$ws.Range('D3').Value = '1.679.69'
$ws.Range('E3').Value = '  -5.41%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.90'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5140'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -12.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.004'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2669'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '21.99'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06369'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07382'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.88%  '
$ws.Range('D12').Value = '1.674.38'
$ws.Range('E12').Value = '  -6.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.574'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5782'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.78%  '
$ws.Range('D15').Value = '1.907.77'
$ws.Range('E15').Value = '  -5.47%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008637'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.14'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -12.83%  '
$ws.Range('D18').Value = '26.367.80'
$ws.Range('E18').Value = '  -7.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.026'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.36%  '
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.90'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '187.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -9.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.246'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.004'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.67'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.570'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1183'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.77'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.07%  '
$ws.Range('E29').Value = '  -5.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05836'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.330'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.523'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.512'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.665'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.007'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.5999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.358'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.668'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.51%  '
$ws.Range('D39').Value = '1.101.33'
$ws.Range('E39').Value = '  -3.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01612'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.894'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.44%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8622'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.004'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.60'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('D45').Value = '1.826.77'
$ws.Range('E45').Value = '  -5.36%  '
$ws.Range('E46').Value = '  +4.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.53'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.003'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.059'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.94%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05221'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.62%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4312'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.35%  '
